$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column B ("Dimension", "Total", "in-gamut")
$ws.Range("B1:D1").EntireColumn.Insert()

# ---- Row 1 (headers) ----
$ws.Range("B1").Value = "Dimension"
$ws.Range("C1").Value = "Total"
$ws.Range("D1").Value = "in-gamut"

# ---- Row 2 (Principal Component Coordinates) ----
$ws.Range("B2").Value = "3"
$ws.Range("C2").Value = "195"
$ws.Range("D2").Value = "167"
$ws.Range("E2").Value = "0.12782"
$ws.Range("F2").Value = "7.73038"
$ws.Range("G2").Value = "0.11218"
$ws.Range("H2").Value = "0.98709"
$ws.Range("I2").Value = "0.00173"
$ws.Range("J2").Value = "2.13649"
$ws.Range("K2").Value = "0.01825"
$ws.Range("L2").Value = "86.38423"
$ws.Range("M2").Value = "0.00143"
$ws.Range("N2").Value = "0.8256"
$ws.Range("O2").Value = "0.72739"
$ws.Range("P2").Value = "1.0"

# ---- Row 3 (XYZ) ----
$ws.Range("B3").Value = "3"
$ws.Range("C3").Value = "195"
$ws.Range("D3").Value = "130"
$ws.Range("E3").Value = "0.05388"
$ws.Range("F3").Value = "3.4832"
$ws.Range("G3").Value = "0.05029"
$ws.Range("H3").Value = "0.99908"
$ws.Range("I3").Value = "0.00266"
$ws.Range("J3").Value = "0.52166"
$ws.Range("K3").Value = "0.36752"
$ws.Range("L3").Value = "28.32849"
$ws.Range("M3").Value = "0.00611"
$ws.Range("N3").Value = "0.32387"
$ws.Range("O3").Value = "0.98862"
$ws.Range("P3").Value = "1.0"

# Resize columns to fit new content (mirrors Excel's bestFit behaviour)
$ws.Range("A1:P3").EntireColumn.AutoFit()

# Nudge the saved window position (matches author's workbook view change)
$wb.Windows.Item(1).Left = 4800
